$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix A2 value: strip the leading non-breaking space from " tomsmith" -> "tomsmith"
$ws.Range("A2").Value = "tomsmith"

# Ensure B2 keeps its password value
$ws.Range("B2").Value = "SuperSecretPassword!"

# Widen column B to fit content (multiple window handling)
# NOTE: Excel's ColumnWidth property is expressed in characters of the Normal
# style font and gets rounded to whole pixels before being persisted as the
# OOXML <col width="..."> attribute, so a small adjustment is needed here to
# land on the target width of exactly 24.
$ws.Columns(2).ColumnWidth = 23.17

# Move the active selection to B8
$ws.Range("B8").Select()
